$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 47
$ws1.Range("F3").Value = 111
$ws1.Range("F4").Value = 1554
$ws1.Range("F5").Value = 246
$ws1.Range("F6").Value = 56
$ws1.Range("F7").Value = 1157
$ws1.Range("F8").Value = 10136
$ws1.Range("F11").Value = 254
$ws1.Range("F13").Value = 384
$ws1.Range("F14").Value = 7010
$ws1.Range("F15").Value = 1094
$ws1.Range("F16").Value = 656
$ws1.Range("F18").Value = 224

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 554

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 47
$ws4.Range("F3").Value = 111
$ws4.Range("F4").Value = 1554
$ws4.Range("F5").Value = 246
$ws4.Range("F7").Value = 56
$ws4.Range("F8").Value = 1157
$ws4.Range("F9").Value = 554
$ws4.Range("F11").Value = 10136
$ws4.Range("F14").Value = 254
$ws4.Range("F16").Value = 384
$ws4.Range("F17").Value = 7010
$ws4.Range("F18").Value = 1094
$ws4.Range("F19").Value = 656
$ws4.Range("F21").Value = 224
